# Regenerate save_data column "K" (column G) values: the sheet previously
# stored a "Strike#" style count; re-derive/rewrite it with the recalculated
# K values (std/mean derived s_vals) per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$newK = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 1
    7  = 0
    8  = 2
    9  = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 0
    14 = 2
    15 = 0
    16 = 0
    17 = 3
    18 = 1
    19 = 3
    20 = 1
    21 = 0
    22 = 3
    23 = 1
    24 = 1
    25 = 3
    26 = 0
    27 = 2
    28 = 0
    29 = 0
    30 = 3
    31 = 1
    32 = 1
    33 = 1
    34 = 3
    35 = 1
    36 = 0
    37 = 2
    38 = 1
    39 = 1
    40 = 1
    41 = 2
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 3
    47 = 1
    48 = 3
    49 = 1
    50 = 2
    51 = 1
    52 = 2
    53 = 2
    54 = 2
    55 = 2
    56 = 0
    57 = 2
    58 = 1
    59 = 1
    60 = 1
    61 = 0
    62 = 1
    63 = 2
    64 = 1
    65 = 2
    66 = 3
    67 = 2
    68 = 1
    69 = 1
    70 = 1
    71 = 2
    72 = 1
    73 = 1
    74 = 1
    76 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
